$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'309.53"
$ws.Range("E2").Value = "'-0.34%"
$ws.Range("D3").Value = "'37.26"
$ws.Range("E3").Value = "'-0.33%"
$ws.Range("D4").Value = "'5.124"
$ws.Range("E4").Value = "'0.52%"
$ws.Range("D5").Value = "'0.07854"
$ws.Range("E5").Value = "'0.80%"
$ws.Range("B6").Value = "KuCoinToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D6").Value = "'8.268"
$ws.Range("E6").Value = "'0.71%"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "'1.880"
$ws.Range("E7").Value = "'-0.08%"
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").Value = "'2.986"
$ws.Range("E8").Value = "'-2.15%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9213"
$ws.Range("E9").Value = "'-0.51%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1067"
$ws.Range("E10").Value = "'-8.61%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1902"
$ws.Range("E11").Value = "'-0.22%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.08919"
$ws.Range("E12").Value = "'-5.06%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03315"
$ws.Range("E13").Value = "'-3.21%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09585"
$ws.Range("E14").Value = "'-0.55%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001382"
$ws.Range("E15").Value = "'0.86%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005906"
$ws.Range("E16").Value = "'1.45%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.393"
$ws.Range("E17").Value = "'-4.01%"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "'4.404"
$ws.Range("E18").Value = "'-0.30%"
$ws.Range("D19").Value = "'0.3449"
$ws.Range("E19").Value = "'1.43%"
$ws.Range("D20").Value = "'6.363"
$ws.Range("E20").Value = "'20.76%"
$ws.Range("E21").Value = "'3.59%"
$ws.Range("D22").Value = "'0.2408"
$ws.Range("E22").Value = "'-7.06%"
$ws.Range("D23").Value = "'0.04358"
$ws.Range("E23").Value = "'0.50%"
$ws.Range("D24").Value = "'0.001198"
$ws.Range("E24").Value = "'0.11%"
$ws.Range("D25").Value = "'0.004275"
$ws.Range("E25").Value = "'0.62%"
$ws.Range("D26").Value = "'0.0001401"
$ws.Range("E26").Value = "'7.72%"
$ws.Range("D27").Value = "'0.0002901"
$ws.Range("D39").Value = "'0.02172"
$ws.Range("E39").Value = "'5.41%"
$ws.Range("D40").Value = "'0.05024"
$ws.Range("E40").Value = "'-0.63%"
$ws.Range("E41").Value = "'-1.07%"
$ws.Range("D42").Value = "'0.1353"
$ws.Range("E42").Value = "'0.63%"
$ws.Range("D43").Value = "'0.008524"
$ws.Range("E43").Value = "'-6.65%"
$ws.Range("D44").Value = "'0.002008"
$ws.Range("E44").Value = "'-2.06%"
$ws.Range("D45").Value = "'0.008104"
$ws.Range("E45").Value = "'-5.98%"
$ws.Range("D46").Value = "'0.00006556"
$ws.Range("E46").Value = "'-2.26%"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'0.00%"
$ws.Range("E48").Value = "'11.81%"
$ws.Range("E49").Value = "'-16.52%"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'0.00%"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'0.00%"
